# Applies updated market-board price/profit figures to the Leve profit
# tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR), as refreshed by the
# scheduled data-pull runner. Columns H-N hold numeric snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: H9=81.25, I9=70.454544, J9=200, K9=70.454544, L9=200, M9=98.545456, N9=-538
$ws.Cells.Item(9, 8).Value = 81.25
$ws.Cells.Item(9, 9).Value = 70.454544
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 70.454544
$ws.Cells.Item(9, 12).Value = 200
$ws.Cells.Item(9, 13).Value = 98.545456
$ws.Cells.Item(9, 14).Value = -538

# Row 28: H28=349.14285, I28=799.5, J28=169, K28=799.5, L28=169, M28=-314.5, N28=-1139
$ws.Cells.Item(28, 8).Value = 349.14285
$ws.Cells.Item(28, 9).Value = 799.5
$ws.Cells.Item(28, 10).Value = 169
$ws.Cells.Item(28, 11).Value = 799.5
$ws.Cells.Item(28, 12).Value = 169
$ws.Cells.Item(28, 13).Value = -314.5
$ws.Cells.Item(28, 14).Value = -1139

# Row 107: H107=214.6923, I107=159.54546, K107=159.54546, M107=1760.45454
$ws.Cells.Item(107, 8).Value = 214.6923
$ws.Cells.Item(107, 9).Value = 159.54546
$ws.Cells.Item(107, 11).Value = 159.54546
$ws.Cells.Item(107, 13).Value = 1760.45454

# Row 116: H116=1788.4, I116=1726.2858, J116=1933.3334, K116=1726.2858, L116=1933.3334, M116=1715.7142, N116=-8817.3334
$ws.Cells.Item(116, 8).Value = 1788.4
$ws.Cells.Item(116, 9).Value = 1726.2858
$ws.Cells.Item(116, 10).Value = 1933.3334
$ws.Cells.Item(116, 11).Value = 1726.2858
$ws.Cells.Item(116, 12).Value = 1933.3334
$ws.Cells.Item(116, 13).Value = 1715.7142
$ws.Cells.Item(116, 14).Value = -8817.3334

$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2=1231.1333, I2=1080.1, J2=1533.2, K2=1080.1, L2=1533.2, M2=-967.0999999999999, N2=-1759.2
$ws.Cells.Item(2, 8).Value = 1231.1333
$ws.Cells.Item(2, 9).Value = 1080.1
$ws.Cells.Item(2, 10).Value = 1533.2
$ws.Cells.Item(2, 11).Value = 1080.1
$ws.Cells.Item(2, 12).Value = 1533.2
$ws.Cells.Item(2, 13).Value = -967.0999999999999
$ws.Cells.Item(2, 14).Value = -1759.2

# Row 32: H32=21169.125, I32=22151.365, J32=8400, K32=22151.365, L32=8400, M32=-21864.365, N32=-8974
$ws.Cells.Item(32, 8).Value = 21169.125
$ws.Cells.Item(32, 9).Value = 22151.365
$ws.Cells.Item(32, 10).Value = 8400
$ws.Cells.Item(32, 11).Value = 22151.365
$ws.Cells.Item(32, 12).Value = 8400
$ws.Cells.Item(32, 13).Value = -21864.365
$ws.Cells.Item(32, 14).Value = -8974

# Row 110: H110=1347.2084, I110=1145, K110=1145, M110=900
$ws.Cells.Item(110, 8).Value = 1347.2084
$ws.Cells.Item(110, 9).Value = 1145
$ws.Cells.Item(110, 11).Value = 1145
$ws.Cells.Item(110, 13).Value = 900

# Row 116: H116=1231.1333, I116=1080.1, J116=1533.2, K116=1080.1, L116=1533.2, M116=1213.9, N116=-6121.2
$ws.Cells.Item(116, 8).Value = 1231.1333
$ws.Cells.Item(116, 9).Value = 1080.1
$ws.Cells.Item(116, 10).Value = 1533.2
$ws.Cells.Item(116, 11).Value = 1080.1
$ws.Cells.Item(116, 12).Value = 1533.2
$ws.Cells.Item(116, 13).Value = 1213.9
$ws.Cells.Item(116, 14).Value = -6121.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3=1231.1333, I3=1080.1, J3=1533.2, K3=1080.1, L3=1533.2, M3=-966.0999999999999, N3=-1761.2
$ws.Cells.Item(3, 8).Value = 1231.1333
$ws.Cells.Item(3, 9).Value = 1080.1
$ws.Cells.Item(3, 10).Value = 1533.2
$ws.Cells.Item(3, 11).Value = 1080.1
$ws.Cells.Item(3, 12).Value = 1533.2
$ws.Cells.Item(3, 13).Value = -966.0999999999999
$ws.Cells.Item(3, 14).Value = -1761.2

# Row 86: H86=1945.4286, I86=1270.3334, J86=2451.75, K86=1270.3334, L86=2451.75, M86=-147.3334, N86=-4697.75
$ws.Cells.Item(86, 8).Value = 1945.4286
$ws.Cells.Item(86, 9).Value = 1270.3334
$ws.Cells.Item(86, 10).Value = 2451.75
$ws.Cells.Item(86, 11).Value = 1270.3334
$ws.Cells.Item(86, 12).Value = 2451.75
$ws.Cells.Item(86, 13).Value = -147.3334
$ws.Cells.Item(86, 14).Value = -4697.75

# Row 89: H89=1945.4286, I89=1270.3334, J89=2451.75, K89=6351.666999999999, L89=12258.75, M89=-735.6669999999995, N89=-23490.75
$ws.Cells.Item(89, 8).Value = 1945.4286
$ws.Cells.Item(89, 9).Value = 1270.3334
$ws.Cells.Item(89, 10).Value = 2451.75
$ws.Cells.Item(89, 11).Value = 6351.666999999999
$ws.Cells.Item(89, 12).Value = 12258.75
$ws.Cells.Item(89, 13).Value = -735.6669999999995
$ws.Cells.Item(89, 14).Value = -23490.75

# Row 107: H107=1258.0416, I107=763.2727, J107=1676.6923, K107=763.2727, L107=1676.6923, M107=1156.7273, N107=-5516.6923
$ws.Cells.Item(107, 8).Value = 1258.0416
$ws.Cells.Item(107, 9).Value = 763.2727
$ws.Cells.Item(107, 10).Value = 1676.6923
$ws.Cells.Item(107, 11).Value = 763.2727
$ws.Cells.Item(107, 12).Value = 1676.6923
$ws.Cells.Item(107, 13).Value = 1156.7273
$ws.Cells.Item(107, 14).Value = -5516.6923

$ws = $wb.Worksheets.Item("CRP")
# Row 16: H16=1737.125, I16=1079.8, J16=2832.6667, K16=1079.8, L16=2832.6667, M16=-792.8, N16=-3406.6667
$ws.Cells.Item(16, 8).Value = 1737.125
$ws.Cells.Item(16, 9).Value = 1079.8
$ws.Cells.Item(16, 10).Value = 2832.6667
$ws.Cells.Item(16, 11).Value = 1079.8
$ws.Cells.Item(16, 12).Value = 2832.6667
$ws.Cells.Item(16, 13).Value = -792.8
$ws.Cells.Item(16, 14).Value = -3406.6667

# Row 86: H86=1916.375, I86=2392, K86=2392, M86=-1269
$ws.Cells.Item(86, 8).Value = 1916.375
$ws.Cells.Item(86, 9).Value = 2392
$ws.Cells.Item(86, 11).Value = 2392
$ws.Cells.Item(86, 13).Value = -1269

# Row 89: H89=1916.375, I89=2392, K89=11960, M89=-6344
$ws.Cells.Item(89, 8).Value = 1916.375
$ws.Cells.Item(89, 9).Value = 2392
$ws.Cells.Item(89, 11).Value = 11960
$ws.Cells.Item(89, 13).Value = -6344

# Row 107: H107=491.93103, I107=490.5, J107=504.33334, K107=490.5, L107=504.33334, M107=1429.5, N107=-4344.33334
$ws.Cells.Item(107, 8).Value = 491.93103
$ws.Cells.Item(107, 9).Value = 490.5
$ws.Cells.Item(107, 10).Value = 504.33334
$ws.Cells.Item(107, 11).Value = 490.5
$ws.Cells.Item(107, 12).Value = 504.33334
$ws.Cells.Item(107, 13).Value = 1429.5
$ws.Cells.Item(107, 14).Value = -4344.33334

# Row 113: H113=1737.125, I113=1079.8, J113=2832.6667, K113=1079.8, L113=2832.6667, M113=1090.2, N113=-7172.6667
$ws.Cells.Item(113, 8).Value = 1737.125
$ws.Cells.Item(113, 9).Value = 1079.8
$ws.Cells.Item(113, 10).Value = 2832.6667
$ws.Cells.Item(113, 11).Value = 1079.8
$ws.Cells.Item(113, 12).Value = 2832.6667
$ws.Cells.Item(113, 13).Value = 1090.2
$ws.Cells.Item(113, 14).Value = -7172.6667

# Row 132: H132=3122.7144, I132=1968.3334, J132=3988.5, K132=5905.0002, L132=11965.5, M132=-3375.0002, N132=-17025.5
$ws.Cells.Item(132, 8).Value = 3122.7144
$ws.Cells.Item(132, 9).Value = 1968.3334
$ws.Cells.Item(132, 10).Value = 3988.5
$ws.Cells.Item(132, 11).Value = 5905.0002
$ws.Cells.Item(132, 12).Value = 11965.5
$ws.Cells.Item(132, 13).Value = -3375.0002
$ws.Cells.Item(132, 14).Value = -17025.5

$ws = $wb.Worksheets.Item("CUL")
# Row 68: H68=1025.3, I68=908.4483, J68=1333.3636, K68=2725.3449, L68=4000.0908, M68=-1914.3449, N68=-5622.0908
$ws.Cells.Item(68, 8).Value = 1025.3
$ws.Cells.Item(68, 9).Value = 908.4483
$ws.Cells.Item(68, 10).Value = 1333.3636
$ws.Cells.Item(68, 11).Value = 2725.3449
$ws.Cells.Item(68, 12).Value = 4000.0908
$ws.Cells.Item(68, 13).Value = -1914.3449
$ws.Cells.Item(68, 14).Value = -5622.0908

# Row 71: H71=1025.3, I71=908.4483, J71=1333.3636, K71=8176.0347, L71=12000.2724, M71=-4120.0347, N71=-20112.2724
$ws.Cells.Item(71, 8).Value = 1025.3
$ws.Cells.Item(71, 9).Value = 908.4483
$ws.Cells.Item(71, 10).Value = 1333.3636
$ws.Cells.Item(71, 11).Value = 8176.0347
$ws.Cells.Item(71, 12).Value = 12000.2724
$ws.Cells.Item(71, 13).Value = -4120.0347
$ws.Cells.Item(71, 14).Value = -20112.2724

# Row 122: H122=811.43475, I122=457, J122=1476, K122=4113, L122=13284, M122=-1663, N122=-18184
$ws.Cells.Item(122, 8).Value = 811.43475
$ws.Cells.Item(122, 9).Value = 457
$ws.Cells.Item(122, 10).Value = 1476
$ws.Cells.Item(122, 11).Value = 4113
$ws.Cells.Item(122, 12).Value = 13284
$ws.Cells.Item(122, 13).Value = -1663
$ws.Cells.Item(122, 14).Value = -18184

$ws = $wb.Worksheets.Item("GSM")
# Row 102: H102=1058.3043, I102=903, J102=1177.7693, K102=903, L102=1177.7693, M102=719, N102=-4421.7693
$ws.Cells.Item(102, 8).Value = 1058.3043
$ws.Cells.Item(102, 9).Value = 903
$ws.Cells.Item(102, 10).Value = 1177.7693
$ws.Cells.Item(102, 11).Value = 903
$ws.Cells.Item(102, 12).Value = 1177.7693
$ws.Cells.Item(102, 13).Value = 719
$ws.Cells.Item(102, 14).Value = -4421.7693

# Row 107: H107=611.7368, I107=547.6923, K107=547.6923, M107=1372.3077
$ws.Cells.Item(107, 8).Value = 611.7368
$ws.Cells.Item(107, 9).Value = 547.6923
$ws.Cells.Item(107, 11).Value = 547.6923
$ws.Cells.Item(107, 13).Value = 1372.3077

# Row 113: H113=11906977, I113=41671500, J113=1167.5333, K113=41671500, L113=1167.5333, M113=-41669330, N113=-5507.5333
$ws.Cells.Item(113, 8).Value = 11906977
$ws.Cells.Item(113, 9).Value = 41671500
$ws.Cells.Item(113, 10).Value = 1167.5333
$ws.Cells.Item(113, 11).Value = 41671500
$ws.Cells.Item(113, 12).Value = 1167.5333
$ws.Cells.Item(113, 13).Value = -41669330
$ws.Cells.Item(113, 14).Value = -5507.5333

$ws = $wb.Worksheets.Item("WVR")
# Row 107: H107=1241.8182, I107=1370.2222, K107=4110.6666, M107=-2190.6666
$ws.Cells.Item(107, 8).Value = 1241.8182
$ws.Cells.Item(107, 9).Value = 1370.2222
$ws.Cells.Item(107, 11).Value = 4110.6666
$ws.Cells.Item(107, 13).Value = -2190.6666

# Row 113: H113=486.81818, I113=333.66666, J113=670.6, K113=1000.99998, L113=2011.8, M113=1169.00002, N113=-6351.8
$ws.Cells.Item(113, 8).Value = 486.81818
$ws.Cells.Item(113, 9).Value = 333.66666
$ws.Cells.Item(113, 10).Value = 670.6
$ws.Cells.Item(113, 11).Value = 1000.99998
$ws.Cells.Item(113, 12).Value = 2011.8
$ws.Cells.Item(113, 13).Value = 1169.00002
$ws.Cells.Item(113, 14).Value = -6351.8

